{"js": "// Replace the multiplication problems in the practice-table cells with a\n// new batch of problems. Each \"old\" value below occurs exactly once in\n// the document, so a plain search + full-text replace on each match is\n// safe and keeps the original run formatting (font, size, etc.) intact.\nconst replacements = [\n  [\"25\u00d723=\", \"86\u00d748=\"],\n  [\"65\u00d787=\", \"25\u00d783=\"],\n  [\"54\u00d742=\", \"86\u00d754=\"],\n  [\"58\u00d794=\", \"91\u00d766=\"],\n  [\"57\u00d773=\", \"84\u00d724=\"],\n  [\"41\u00d747=\", \"18\u00d763=\"],\n  [\"29\u00d752=\", \"81\u00d790=\"],\n  [\"95\u00d721=\", \"23\u00d756=\"],\n  [\"86\u00d792=\", \"68\u00d765=\"],\n  [\"84\u00d774=\", \"68\u00d741=\"],\n  [\"65\u00d773=\", \"88\u00d793=\"],\n  [\"27\u00d722=\", \"80\u00d786=\"],\n  [\"33\u00d789=\", \"53\u00d741=\"],\n  [\"13\u00d750=\", \"51\u00d758=\"],\n  [\"87\u00d735=\", \"26\u00d783=\"],\n  [\"38\u00d778=\", \"75\u00d780=\"],\n  [\"18\u00d716=\", \"47\u00d797=\"],\n  [\"71\u00d782=\", \"28\u00d733=\"],\n  [\"89\u00d740=\", \"29\u00d777=\"],\n  [\"35\u00d742=\", \"32\u00d728=\"],\n  [\"65\u00d780=\", \"51\u00d711=\"],\n  [\"12\u00d735=\", \"83\u00d749=\"],\n  [\"33\u00d751=\", \"20\u00d774=\"],\n  [\"76\u00d771=\", \"93\u00d790=\"],\n  [\"62\u00d726=\", \"96\u00d798=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication problems in the practice-table cells with a\n# new batch of problems. Each \"Old\" value below occurs exactly once in\n# the document, so Find/Replace (wdReplaceOne) on the whole-document\n# range for each pair is safe and preserves the original run formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"25\u00d723=\"; New = \"86\u00d748=\" },\n    @{ Old = \"65\u00d787=\"; New = \"25\u00d783=\" },\n    @{ Old = \"54\u00d742=\"; New = \"86\u00d754=\" },\n    @{ Old = \"58\u00d794=\"; New = \"91\u00d766=\" },\n    @{ Old = \"57\u00d773=\"; New = \"84\u00d724=\" },\n    @{ Old = \"41\u00d747=\"; New = \"18\u00d763=\" },\n    @{ Old = \"29\u00d752=\"; New = \"81\u00d790=\" },\n    @{ Old = \"95\u00d721=\"; New = \"23\u00d756=\" },\n    @{ Old = \"86\u00d792=\"; New = \"68\u00d765=\" },\n    @{ Old = \"84\u00d774=\"; New = \"68\u00d741=\" },\n    @{ Old = \"65\u00d773=\"; New = \"88\u00d793=\" },\n    @{ Old = \"27\u00d722=\"; New = \"80\u00d786=\" },\n    @{ Old = \"33\u00d789=\"; New = \"53\u00d741=\" },\n    @{ Old = \"13\u00d750=\"; New = \"51\u00d758=\" },\n    @{ Old = \"87\u00d735=\"; New = \"26\u00d783=\" },\n    @{ Old = \"38\u00d778=\"; New = \"75\u00d780=\" },\n    @{ Old = \"18\u00d716=\"; New = \"47\u00d797=\" },\n    @{ Old = \"71\u00d782=\"; New = \"28\u00d733=\" },\n    @{ Old = \"89\u00d740=\"; New = \"29\u00d777=\" },\n    @{ Old = \"35\u00d742=\"; New = \"32\u00d728=\" },\n    @{ Old = \"65\u00d780=\"; New = \"51\u00d711=\" },\n    @{ Old = \"12\u00d735=\"; New = \"83\u00d749=\" },\n    @{ Old = \"33\u00d751=\"; New = \"20\u00d774=\" },\n    @{ Old = \"76\u00d771=\"; New = \"93\u00d790=\" },\n    @{ Old = \"62\u00d726=\"; New = \"96\u00d798=\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
